$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh - update D (Price) and E (Volume 1h) columns
# per-row, matching upstream scrape snapshot. Numeric-looking text in column D
# is written via a temporary Text number format so Excel keeps it as a literal
# string (preserving formats like "1.001" / "0.06063" / trailing zeros) instead
# of silently coercing it to a float; ClearFormats() afterwards restores the
# cell to its original unstyled state so only the value itself changes.

$ws.Range('D2').Value = '25.666.43'
$ws.Range('E2').Value = '  -3.25%  '
$ws.Range('D3').Value = '1.740.77'
$ws.Range('E3').Value = '  -5.63%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.86'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -10.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4897'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -8.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.23'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -8.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2576'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -17.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06063'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -12.10%  '
$ws.Range('D11').Value = '1.742.53'
$ws.Range('E11').Value = '  -5.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06832'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -12.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.75'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -20.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.447'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -11.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '75.53'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -16.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.5625'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -26.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '25.696.63'
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.40'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -18.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006574'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -17.33%  '
$ws.Range('D22').Value = '1.962.27'
$ws.Range('E22').Value = '  -5.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.026'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -13.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.892'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -15.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.995'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -17.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '137.27'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.448'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -14.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.821'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -16.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.67'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -14.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '101.15'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -9.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.699'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -13.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07968'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -9.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.403'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -17.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04401'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -8.99%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.612'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -11.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9733'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -14.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5916'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -19.45%  '
$ws.Range('E39').Value = '  -14.51%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '103.01'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -5.08%  '
$ws.Range('E42').Value = '  -12.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.831'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -21.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.146'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3730'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -22.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7203'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -20.33%  '
$ws.Range('E47').Value = '  -10.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1079'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -13.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.90'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -14.54%  '
$ws.Range('E50').Value = '  -14.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.781'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -24.47%  '
